# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# The source data rows for a handful of matches were in the wrong order
# relative to their match id (column B). This re-sorts those rows by
# swapping everything except column A (the running index) between the
# affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")

    $v1 = $rng1.Value()
    $v2 = $rng2.Value()

    $rng1.Value = $v2
    $rng2.Value = $v1
}

# Simple two-row swaps
Swap-Rows $ws 35 36
Swap-Rows $ws 62 63
Swap-Rows $ws 75 76
Swap-Rows $ws 77 78
Swap-Rows $ws 226 227

# Three-row rotation: row 228 <- row 231 <- row 230 <- row 228
# (row 229 in between is untouched)
$rng228 = $ws.Range("B228:AD228")
$rng230 = $ws.Range("B230:AD230")
$rng231 = $ws.Range("B231:AD231")

$v228 = $rng228.Value()
$v230 = $rng230.Value()
$v231 = $rng231.Value()

$rng228.Value = $v231
$rng230.Value = $v228
$rng231.Value = $v230
